$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.501.84'
$ws.Range("E2").Value = '  +9.47%  '

$ws.Range("D3").Value = '1.596.10'
$ws.Range("E3").Value = '  +8.37%  '

$ws.Range("D4").Value = '0.9959'
$ws.Range("E4").Value = '  -1.08%  '

$ws.Range("D5").Value = '0.9896'
$ws.Range("E5").Value = '  +2.88%  '

$ws.Range("D6").Value = '299.44'
$ws.Range("E6").Value = '  +8.00%  '

$ws.Range("D7").Value = '0.3612'
$ws.Range("E7").Value = '  +0.60%  '

$ws.Range("D8").Value = '0.3348'
$ws.Range("E8").Value = '  +8.68%  '

$ws.Range("D9").Value = '''41.10'
$ws.Range("E9").Value = '  +4.27%  '

$ws.Range("D10").Value = '1.115'
$ws.Range("E10").Value = '  +2.55%  '

$ws.Range("D11").Value = '0.06923'
$ws.Range("E11").Value = '  +4.47%  '

$ws.Range("D12").Value = '0.9924'
$ws.Range("E12").Value = '  -0.89%  '

$ws.Range("D13").Value = '19.39'
$ws.Range("E13").Value = '  +6.82%  '

$ws.Range("D14").Value = '5.791'
$ws.Range("E14").Value = '  +5.96%  '

$ws.Range("D15").Value = '6.486'
$ws.Range("E15").Value = '  +5.15%  '

$ws.Range("D16").Value = '0.9887'
$ws.Range("E16").Value = '  +2.80%  '

$ws.Range("D17").Value = '0.00001063'
$ws.Range("E17").Value = '  +3.95%  '

$ws.Range("D18").Value = '1.591.66'
$ws.Range("E18").Value = '  +8.07%  '

$ws.Range("D19").Value = '0.06581'
$ws.Range("E19").Value = '  +10.04%  '

$ws.Range("D20").Value = '76.29'
$ws.Range("E20").Value = '  +10.59%  '

$ws.Range("D21").Value = '15.85'
$ws.Range("E21").Value = '  +8.85%  '

$ws.Range("D22").Value = '''5.900'
$ws.Range("E22").Value = '  +7.50%  '

$ws.Range("D23").Value = '11.47'
$ws.Range("E23").Value = '  +1.97%  '

$ws.Range("D24").Value = '22.432.40'
$ws.Range("E24").Value = '  +9.15%  '

$ws.Range("E25").Value = '  +4.49%  '

$ws.Range("D26").Value = '''2.490'
$ws.Range("E26").Value = '  +18.32%  '

$ws.Range("D27").Value = '149.22'
$ws.Range("E27").Value = '  +3.41%  '

$ws.Range("D28").Value = '19.13'
$ws.Range("E28").Value = '  +11.68%  '

$ws.Range("D29").Value = '1.761.17'
$ws.Range("E29").Value = '  +7.80%  '

$ws.Range("D30").Value = '122.69'
$ws.Range("E30").Value = '  +7.81%  '

$ws.Range("D31").Value = '3.931'
$ws.Range("E31").Value = '  +1.25%  '

$ws.Range("D32").Value = '5.844'
$ws.Range("E32").Value = '  +18.42%  '

$ws.Range("D33").Value = '0.9217'
$ws.Range("E33").Value = '  +14.63%  '

$ws.Range("D34").Value = '0.08089'
$ws.Range("E34").Value = '  +1.20%  '

$ws.Range("D35").Value = '1.626'
$ws.Range("E35").Value = '  +10.60%  '

$ws.Range("D36").Value = '11.75'
$ws.Range("E36").Value = '  +13.33%  '

$ws.Range("D37").Value = '1.238'
$ws.Range("E37").Value = '  -0.67%  '

$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").Value = '8.416'
$ws.Range("E38").Value = '  +14.19%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '5.043'
$ws.Range("E39").Value = '  +7.08%  '

$ws.Range("D40").Value = '0.05977'
$ws.Range("E40").Value = '  +3.42%  '

$ws.Range("D41").Value = '0.02178'
$ws.Range("E41").Value = '  +6.42%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1973'
$ws.Range("E42").Value = '  +5.16%  '

$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").Value = '0.9897'
$ws.Range("E43").Value = '  +2.81%  '

$ws.Range("D44").Value = '0.5764'
$ws.Range("E44").Value = '  +9.54%  '

$ws.Range("D45").Value = '3.756'
$ws.Range("E45").Value = '  +6.74%  '

$ws.Range("D46").Value = '12.68'
$ws.Range("E46").Value = '  +4.77%  '

$ws.Range("D47").Value = '0.5591'
$ws.Range("E47").Value = '  +7.63%  '

$ws.Range("D48").Value = '''123.30'
$ws.Range("E48").Value = '  +3.57%  '

$ws.Range("D49").Value = '1.935'
$ws.Range("E49").Value = '  +6.98%  '

$ws.Range("D50").Value = '0.06779'
$ws.Range("E50").Value = '  +5.19%  '

$ws.Range("D51").Value = '72.34'
$ws.Range("E51").Value = '  +7.80%  '

